$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.00037110940320417285

$ws.Range("A3").Value = 0.00009978991147363558
$ws.Range("C3").Value = 42.20338821411133
$ws.Range("D3").Value = 18.982187271118164

$ws.Range("A4").Value = 0.00009865831088973209
$ws.Range("C4").Value = 38.64406967163086
$ws.Range("D4").Value = 17.38983154296875

$ws.Range("A5").Value = 0.000057155593822244555
$ws.Range("C5").Value = 49.830509185791016
$ws.Range("D5").Value = 35.906951904296875

$ws.Range("A6").Value = 0.00002646152461238671
$ws.Range("C6").Value = 42.20338821411133
$ws.Range("D6").Value = 18.964170455932617

$ws.Range("A7").Value = 0.000025632543838582933
$ws.Range("C7").Value = 46.1016960144043
$ws.Range("D7").Value = 20.745763778686523

$ws.Range("A8").Value = 0.000024724577087908983
$ws.Range("C8").Value = 42.37288284301758
$ws.Range("D8").Value = 19.06779670715332

$ws.Range("A9").Value = 0.000023012544261291623
$ws.Range("C9").Value = 38.64406967163086
$ws.Range("D9").Value = 17.38983154296875

$ws.Range("A10").Value = 0.000005139661425346276
$ws.Range("C10").Value = 38.64406967163086
$ws.Range("D10").Value = 17.38983154296875

$ws.Range("A11").Value = 0.000004932881438435288
$ws.Range("C11").Value = 46.1016960144043
$ws.Range("D11").Value = 17.38983154296875

$ws.Range("A12").Value = 0.00000426711858381168
$ws.Range("C12").Value = 4.0677971839904785
$ws.Range("D12").Value = 1.8305089473724365

$ws.Range("A13").Value = 0.0000013347457752388436
$ws.Range("C13").Value = 42.37288284301758
$ws.Range("D13").Value = 18.86345672607422
